# edit.ps1 - PowerPoint COM-interop script
#
# Reproduces:
#  1. The table on slide 6 is switched to a different (built-in) table
#     style, identified by its style GUID.
#  2. The presentation's theme colour scheme is swapped for the one that
#     used to live in the "other" theme part ("Office Theme" <-> "Integral"),
#     i.e. the deck's 12 theme colours are rewritten to the classic default
#     Office palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{6F80C769-604E-4CFD-A005-8FBAC7272E06}")

# --- 2. Swap the theme colour scheme -------------------------------------
# Small helper replicating VBA's RGB() (R + G*256 + B*65536), since this
# host does not expose an "RGB" cmdlet.
function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Colour values that currently live in the "Office Theme" theme part; they
# become the deck's active theme colours (indices follow the standard
# ThemeColorScheme ordering: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
$officeThemeColors = @(
    (RGBVal 0x00 0x00 0x00),
    (RGBVal 0xFF 0xFF 0xFF),
    (RGBVal 0x44 0x54 0x6A),
    (RGBVal 0xE7 0xE6 0xE6),
    (RGBVal 0x5B 0x9B 0xD5),
    (RGBVal 0xED 0x7D 0x31),
    (RGBVal 0xA5 0xA5 0xA5),
    (RGBVal 0xFF 0xC0 0x00),
    (RGBVal 0x44 0x72 0xC4),
    (RGBVal 0x70 0xAD 0x47),
    (RGBVal 0x05 0x63 0xC1),
    (RGBVal 0x95 0x4F 0x72)
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeColors[$i - 1]
}
